# Regenerate the quadratic/linear problem data (per commit message:
# "volver a generar problemas cuadraticos y lineales").
#
# All of the values below are stored in the workbook as TEXT (not numbers) -
# this mirrors how the source file was produced (python -> pandas -> xlsx),
# so every cell we touch gets a text number-format applied first to stop
# Excel from re-interpreting numeric-looking strings ("0.21", "0", ...) as
# real numbers. The number format is reset back to the default afterwards
# so the written workbook doesn't pick up any stray text formatting.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- Restricciones_del_follower ---------------------------------------
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

Set-TextValue $ws.Range("A2") "0.9000000000000004 - 3x + y"
Set-TextValue $ws.Range("B2") "2.0999999999999996"
Set-TextValue $ws.Range("D2") "0.13"
Set-TextValue $ws.Range("E2") "2.8000000000000003"
# F2 keeps its old value ("0") - only its *meaning* shifts because F3 (which
# used to share the same text) changes independently below.

Set-TextValue $ws.Range("A3") "0.47499999999999964 + x - 0.5y"
Set-TextValue $ws.Range("B3") "-4.475"
Set-TextValue $ws.Range("D3") "0.21"
Set-TextValue $ws.Range("E3") "0"
Set-TextValue $ws.Range("F3") "2.8000000000000003"

Set-TextValue $ws.Range("A4") "-7.5 + x + y"
Set-TextValue $ws.Range("B4") "-0.49999999999999956"
Set-TextValue $ws.Range("D4") "0.91"
Set-TextValue $ws.Range("E4") "0"
Set-TextValue $ws.Range("F4") "6.1"

# --- Punto_modificado ---------------------------------------------------
$ws = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws.Range("A2") "1.85"
Set-TextValue $ws.Range("B2") "4.65"

# --- Vector_bf / Vector_BF ------------------------------------------------
# NOTE: worksheet lookup by name is case-insensitive, and this workbook has
# two sheets whose names differ only by case ("Vector_bf" / "Vector_BF"), so
# they must be addressed by their (1-based) tab position instead, or the
# second lookup silently resolves back to the first sheet.
$ws = $wb.Worksheets.Item(5)   # Vector_bf
Set-TextValue $ws.Range("A2") "-5.460000000000001"

$ws = $wb.Worksheets.Item(6)   # Vector_BF
Set-TextValue $ws.Range("A2") "14.7"
Set-TextValue $ws.Range("A3") "-44.0"

Write-Output "edit applied"
